# Auto-generated edit script: update crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.670.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.294.27"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.83"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.17"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.57%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.16%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.81"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.60"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.14%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.83"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.652.00"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.300.93"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.773"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.613.53"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.74"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0892"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.99"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.07"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.86"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.38%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.38"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.64"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.49"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.05"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.98"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.87"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.70"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.94"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.43"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.36"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0683"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.991.57"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0279"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.13"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.22"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.89%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.76"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.91"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.45"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.518.73"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.74"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.01%  "
